$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1078.1875
$ws.Range("I100").Value = 1063.7273
$ws.Range("J100").Value = 1110
$ws.Range("K100").Value = 1063.7273
$ws.Range("L100").Value = 1110
$ws.Range("M100").Value = -522.7273
$ws.Range("N100").Value = -2192
$ws.Range("H137").Value = 1388.1708
$ws.Range("I137").Value = 1248.1875
$ws.Range("J137").Value = 1885.8889
$ws.Range("K137").Value = 3744.5625
$ws.Range("L137").Value = 5657.6667
$ws.Range("M137").Value = -1194.5625
$ws.Range("N137").Value = -10757.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1120.7667
$ws.Range("I2").Value = 853.5294
$ws.Range("K2").Value = 853.5294
$ws.Range("M2").Value = -740.5294
$ws.Range("H32").Value = 3482.0186
$ws.Range("I32").Value = 3008.3
$ws.Range("K32").Value = 3008.3
$ws.Range("M32").Value = -2721.3
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()
$ws.Range("H61").Value = 1339.3103
$ws.Range("I61").Value = 1142.2222
$ws.Range("K61").Value = 1142.2222
$ws.Range("M61").Value = -930.2221999999999
$ws.Range("H74").Value = 1931.4348
$ws.Range("I74").Value = 1564.3684
$ws.Range("J74").Value = 3675
$ws.Range("K74").Value = 1564.3684
$ws.Range("L74").Value = 3675
$ws.Range("M74").Value = -690.3684000000001
$ws.Range("N74").Value = -5423
$ws.Range("H77").Value = 1931.4348
$ws.Range("I77").Value = 1564.3684
$ws.Range("J77").Value = 3675
$ws.Range("K77").Value = 7821.842000000001
$ws.Range("L77").Value = 18375
$ws.Range("M77").Value = -3453.842000000001
$ws.Range("N77").Value = -27111
$ws.Range("H97").Value = 779.82355
$ws.Range("I97").Value = 659
$ws.Range("J97").Value = 1172.5
$ws.Range("K97").Value = 659
$ws.Range("L97").Value = 1172.5
$ws.Range("M97").Value = -163
$ws.Range("N97").Value = -2164.5
$ws.Range("H98").Value = 18000
$ws.Range("J98").Value = 18000
$ws.Range("L98").Value = 18000
$ws.Range("N98").Value = -23990
$ws.Range("H110").Value = 1103.1538
$ws.Range("I110").Value = 849.1818
$ws.Range("J110").Value = 2500
$ws.Range("K110").Value = 849.1818
$ws.Range("L110").Value = 2500
$ws.Range("M110").Value = 1195.8182
$ws.Range("N110").Value = -6590
$ws.Range("H116").Value = 1120.7667
$ws.Range("I116").Value = 853.5294
$ws.Range("K116").Value = 853.5294
$ws.Range("M116").Value = 1440.4706
$ws.Range("H122").Value = 1594.25
$ws.Range("I122").Value = 1142.125
$ws.Range("J122").Value = 2498.5
$ws.Range("K122").Value = 3426.375
$ws.Range("L122").Value = 7495.5
$ws.Range("M122").Value = -976.375
$ws.Range("N122").Value = -12395.5
$ws.Range("H132").Value = 4413.2573
$ws.Range("I132").Value = 4601.3105
$ws.Range("J132").Value = 3504.3333
$ws.Range("K132").Value = 13803.9315
$ws.Range("L132").Value = 10512.9999
$ws.Range("M132").Value = -11273.9315
$ws.Range("N132").Value = -15572.9999
$ws.Range("H136").Value = 1339.3103
$ws.Range("I136").Value = 1142.2222
$ws.Range("K136").Value = 3426.6666
$ws.Range("M136").Value = -876.6665999999996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1120.7667
$ws.Range("I3").Value = 853.5294
$ws.Range("K3").Value = 853.5294
$ws.Range("M3").Value = -739.5294
$ws.Range("H94").Value = 1083.1111
$ws.Range("I94").Value = 862.25
$ws.Range("K94").Value = 862.25
$ws.Range("M94").Value = -411.25
$ws.Range("H95").Value = 22300
$ws.Range("J95").Value = 22300
$ws.Range("L95").Value = 22300
$ws.Range("N95").Value = -27792
$ws.Range("H99").Value = 654.4
$ws.Range("I99").Value = 505.5
$ws.Range("K99").Value = 505.5
$ws.Range("M99").Value = 992.5
$ws.Range("H107").Value = 1535
$ws.Range("I107").Value = 1227.75
$ws.Range("J107").Value = 2149.5
$ws.Range("K107").Value = 1227.75
$ws.Range("L107").Value = 2149.5
$ws.Range("M107").Value = 692.25
$ws.Range("N107").Value = -5989.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 32360.75
$ws.Range("J28").Value = 32360.75
$ws.Range("L28").Value = 32360.75
$ws.Range("N28").Value = -32850.75
$ws.Range("H96").Value = 26899.5
$ws.Range("J96").Value = 26899.5
$ws.Range("L96").Value = 26899.5
$ws.Range("N96").Value = -32391.5
$ws.Range("H99").Value = 2553.8462
$ws.Range("I99").Value = 1925
$ws.Range("J99").Value = 2833.3333
$ws.Range("K99").Value = 1925
$ws.Range("L99").Value = 2833.3333
$ws.Range("M99").Value = -427
$ws.Range("N99").Value = -5829.3333
$ws.Range("H122").Value = 867.5714
$ws.Range("I122").Value = 867.5714
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2602.7142
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -152.7142000000003
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 2553.8462
$ws.Range("I126").Value = 1925
$ws.Range("J126").Value = 2833.3333
$ws.Range("K126").Value = 5775
$ws.Range("L126").Value = 8499.999899999999
$ws.Range("M126").Value = -3305
$ws.Range("N126").Value = -13439.9999
$ws.Range("H127").Value = 40000
$ws.Range("J127").Value = 40000
$ws.Range("L127").Value = 40000
$ws.Range("N127").Value = -49920
$ws.Range("H132").Value = 3229.65
$ws.Range("I132").Value = 2599.6875
$ws.Range("J132").Value = 5749.5
$ws.Range("K132").Value = 7799.0625
$ws.Range("L132").Value = 17248.5
$ws.Range("M132").Value = -5269.0625
$ws.Range("N132").Value = -22308.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 390
$ws.Range("I68").Value = 385
$ws.Range("J68").Value = 400
$ws.Range("K68").Value = 1155
$ws.Range("L68").Value = 1200
$ws.Range("M68").Value = -344
$ws.Range("N68").Value = -2822
$ws.Range("H71").Value = 390
$ws.Range("I71").Value = 385
$ws.Range("J71").Value = 400
$ws.Range("K71").Value = 3465
$ws.Range("L71").Value = 3600
$ws.Range("M71").Value = 591
$ws.Range("N71").Value = -11712
$ws.Range("H115").Value = 4207.2
$ws.Range("I115").Value = 4014
$ws.Range("J115").Value = 4980
$ws.Range("K115").Value = 12042
$ws.Range("L115").Value = 14940
$ws.Range("M115").Value = -10867
$ws.Range("N115").Value = -17290
$ws.Range("H131").Value = 4450.926
$ws.Range("I131").Value = 10508.8
$ws.Range("J131").Value = 887.4706
$ws.Range("K131").Value = 31526.4
$ws.Range("L131").Value = 2662.4118
$ws.Range("M131").Value = -26486.4
$ws.Range("N131").Value = -12742.4118
$ws.Range("H134").Value = 1701.8
$ws.Range("I134").Value = 1552.7
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 4658.1
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = 411.8999999999996
$ws.Range("N134").Value = -16140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 79238.734
$ws.Range("I132").Value = 96848.09
$ws.Range("J132").Value = 5279.4
$ws.Range("K132").Value = 290544.27
$ws.Range("L132").Value = 15838.2
$ws.Range("M132").Value = -288014.27
$ws.Range("N132").Value = -20898.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2537.5
$ws.Range("I7").Value = 1766.6666
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 1766.6666
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -1654.6666
$ws.Range("N7").Value = -3224
$ws.Range("H40").Value = 1971.75
$ws.Range("I40").Value = 1910.5714
$ws.Range("J40").Value = 2400
$ws.Range("K40").Value = 1910.5714
$ws.Range("L40").Value = 2400
$ws.Range("M40").Value = -1774.5714
$ws.Range("N40").Value = -2672
$ws.Range("H61").Value = 1255.1818
$ws.Range("I61").Value = 980.7
$ws.Range("K61").Value = 980.7
$ws.Range("M61").Value = -778.7
$ws.Range("H93").Value = 1352389.8
$ws.Range("I93").Value = 1931464.8
$ws.Range("J93").Value = 1214.6666
$ws.Range("K93").Value = 1931464.8
$ws.Range("L93").Value = 1214.6666
$ws.Range("M93").Value = -1930216.8
$ws.Range("N93").Value = -3710.6666
$ws.Range("H100").Value = 1375
$ws.Range("I100").Value = 1354.1666
$ws.Range("K100").Value = 1354.1666
$ws.Range("M100").Value = -813.1666
$ws.Range("H113").Value = 1255.1818
$ws.Range("I113").Value = 980.7
$ws.Range("K113").Value = 980.7
$ws.Range("M113").Value = 1189.3
$ws.Range("H126").Value = 2537.5
$ws.Range("I126").Value = 1766.6666
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 5299.9998
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -2829.9998
$ws.Range("N126").Value = -13940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 10000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 10000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -10298
$ws.Range("H81").Value = 2506.25
$ws.Range("I81").Value = 2571.5715
$ws.Range("J81").Value = 2455.4443
$ws.Range("K81").Value = 5143.143
$ws.Range("L81").Value = 4910.8886
$ws.Range("M81").Value = -4082.143
$ws.Range("N81").Value = -7032.8886
$ws.Range("H84").Value = 2506.25
$ws.Range("I84").Value = 2571.5715
$ws.Range("J84").Value = 2455.4443
$ws.Range("K84").Value = 25715.715
$ws.Range("L84").Value = 24554.443
$ws.Range("M84").Value = -20411.715
$ws.Range("N84").Value = -35162.443
$ws.Range("H100").Value = 586.6667
$ws.Range("I100").Value = 454.2857
$ws.Range("J100").Value = 772
$ws.Range("K100").Value = 908.5714
$ws.Range("L100").Value = 1544
$ws.Range("M100").Value = -367.5714
$ws.Range("N100").Value = -2626
$ws.Range("H132").Value = 1725.3158
$ws.Range("I132").Value = 1173.9375
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 3521.8125
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -991.8125
$ws.Range("N132").Value = -19058
